# Apply the LinuxForHealth re-branding update to the ReengagementValue
# StructureDefinition workbook.
#
# Sheet "Metadata" holds a simple Property/Value table (columns A/B).
# Sheet "Elements" holds the FHIR element definitions table; the only
# real content change there is that the ele-1/ext-1 constraint text
# moves from the root "Extension" row down to the "Extension.extension"
# row (column "Constraint(s)").

$wb = $excel.ActiveWorkbook

# ---- Sheet: Metadata ----
$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/reengagement-value"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

# ---- Sheet: Elements ----
$elements = $wb.Worksheets.Item("Elements")

$constraintText = "ele-1:All FHIR elements must have a @value or children {hasValue() or (children().count() > id.count())}`next-1:Must have either extensions or value[x], not both {extension.exists() != value.exists()}"

# Row 2 = "Extension" (root element) - Constraint(s) column AI no longer holds the text
$elements.Range("AI2").Value = ""

# Row 4 = "Extension.extension" - Constraint(s) column AI now holds the text
$elements.Range("AI4").Value = $constraintText

# Row 5 = "Extension.url" - Fixed Value column Q mirrors the StructureDefinition URL
$elements.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/reengagement-value"

$wb.Save()
